$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-07 Friday" "2025-03-08 Saturday"

Replace-Text "46×61=2806" "89×17=1513"
Replace-Text "28×26=728" "94×69=6486"
Replace-Text "84×20=1680" "63×31=1953"
Replace-Text "47×25=1175" "35×51=1785"
Replace-Text "63×74=4662" "18×44=792"
Replace-Text "48×14=672" "13×65=845"
Replace-Text "98×60=5880" "54×40=2160"
Replace-Text "18×86=1548" "13×44=572"
Replace-Text "39×28=1092" "67×65=4355"
Replace-Text "63×83=5229" "20×31=620"
Replace-Text "15×38=570" "32×77=2464"
Replace-Text "82×18=1476" "43×73=3139"
Replace-Text "86×84=7224" "84×90=7560"
Replace-Text "45×38=1710" "78×64=4992"
Replace-Text "27×70=1890" "12×57=684"
Replace-Text "23×93=2139" "46×29=1334"
Replace-Text "42×37=1554" "40×48=1920"
Replace-Text "26×45=1170" "13×50=650"
Replace-Text "50×39=1950" "46×16=736"
Replace-Text "60×92=5520" "79×56=4424"
Replace-Text "34×43=1462" "46×88=4048"
Replace-Text "72×20=1440" "46×13=598"
Replace-Text "53×44=2332" "37×68=2516"
Replace-Text "60×97=5820" "73×52=3796"
Replace-Text "21×51=1071" "53×81=4293"
